# Fruta / hortaliza, semanal
# Insert a new weekly record as row 24 (pushing the existing rows 24-44
# down to 25-45), matching the new "Primera" quality Damasco entry for
# Feria Lagunitas de Puerto Montt / Región Metropolitana.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 24:44 down to 25:45, creating a blank row 24.
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with the new record.
$ws.Cells.Item(24, 1).Value = 4
$ws.Cells.Item(24, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(24, 3).Value = "Los Lagos"
$ws.Cells.Item(24, 4).Value = 44894
$ws.Cells.Item(24, 5).Value = 10
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100103
$ws.Cells.Item(24, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(24, 9).Value = 100103003
$ws.Cells.Item(24, 10).Value = "Damasco"
$ws.Cells.Item(24, 11).Value = "Castle Brite"
$ws.Cells.Item(24, 12).Value = "Primera"
$ws.Cells.Item(24, 13).Value = 400
$ws.Cells.Item(24, 14).Value = 23000
$ws.Cells.Item(24, 15).Value = 24000
$ws.Cells.Item(24, 16).Value = 23500
$ws.Cells.Item(24, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(24, 18).Value = "Región Metropolitana"
$ws.Cells.Item(24, 19).Value = 1469
$ws.Cells.Item(24, 20).Value = 16
